$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1986301369863014
$ws.Cells.Item(2, 3).Value = 0.547945205479452
$ws.Cells.Item(2, 10).Value = 0.00684931506849315
$ws.Cells.Item(2, 16).Value = 0.1643835616438356
$ws.Cells.Item(2, 19).Value = 0.0821917808219178
$ws.Cells.Item(3, 3).Value = 0.01818181818181818
$ws.Cells.Item(3, 10).Value = 0.02424242424242424
$ws.Cells.Item(3, 16).Value = 0.6909090909090909
$ws.Cells.Item(3, 19).Value = 0.2666666666666667
$ws.Cells.Item(4, 10).Value = 0.05263157894736842
$ws.Cells.Item(4, 15).Value = 0.02631578947368421
$ws.Cells.Item(4, 16).Value = 0.6842105263157895
$ws.Cells.Item(4, 19).Value = 0.2368421052631579
$ws.Cells.Item(6, 2).Value = 0.04390243902439024
$ws.Cells.Item(6, 4).Value = 0.01951219512195122
$ws.Cells.Item(6, 6).Value = 0.03414634146341464
$ws.Cells.Item(6, 10).Value = 0.2829268292682927
$ws.Cells.Item(6, 15).Value = 0.01463414634146342
$ws.Cells.Item(6, 17).Value = 0.1268292682926829
$ws.Cells.Item(6, 18).Value = 0.08780487804878048
$ws.Cells.Item(6, 19).Value = 0.3902439024390244
$ws.Cells.Item(7, 2).Value = 0.09663865546218488
$ws.Cells.Item(7, 4).Value = 0.02100840336134454
$ws.Cells.Item(7, 6).Value = 0.02100840336134454
$ws.Cells.Item(7, 10).Value = 0.1764705882352941
$ws.Cells.Item(7, 15).Value = 0.02100840336134454
$ws.Cells.Item(7, 17).Value = 0.1302521008403361
$ws.Cells.Item(7, 18).Value = 0.08403361344537816
$ws.Cells.Item(7, 19).Value = 0.4495798319327731
$ws.Cells.Item(8, 2).Value = 0.08723404255319149
$ws.Cells.Item(8, 4).Value = 0.01276595744680851
$ws.Cells.Item(8, 6).Value = 0.0425531914893617
$ws.Cells.Item(8, 10).Value = 0.1404255319148936
$ws.Cells.Item(8, 15).Value = 0.002127659574468085
$ws.Cells.Item(8, 17).Value = 0.1659574468085106
$ws.Cells.Item(8, 18).Value = 0.07659574468085106
$ws.Cells.Item(8, 19).Value = 0.4723404255319149
$ws.Cells.Item(9, 2).Value = 0.08163265306122448
$ws.Cells.Item(9, 4).Value = 0.01530612244897959
$ws.Cells.Item(9, 6).Value = 0.08673469387755102
$ws.Cells.Item(9, 10).Value = 0.1173469387755102
$ws.Cells.Item(9, 15).Value = 0.01020408163265306
$ws.Cells.Item(9, 17).Value = 0.1836734693877551
$ws.Cells.Item(9, 18).Value = 0.1071428571428571
$ws.Cells.Item(9, 19).Value = 0.3979591836734694
$ws.Cells.Item(10, 2).Value = 0.1098221191028616
$ws.Cells.Item(10, 4).Value = 0.01701469450889405
$ws.Cells.Item(10, 5).Value = 0.0007733952049497294
$ws.Cells.Item(10, 6).Value = 0.0711523588553751
$ws.Cells.Item(10, 10).Value = 0.1160092807424594
$ws.Cells.Item(10, 15).Value = 0.01469450889404486
$ws.Cells.Item(10, 17).Value = 0.1832946635730859
$ws.Cells.Item(10, 18).Value = 0.07965970610982212
$ws.Cells.Item(10, 19).Value = 0.4075792730085073
$ws.Cells.Item(11, 7).Value = 0.1428571428571428
$ws.Cells.Item(11, 10).Value = 0.09022556390977443
$ws.Cells.Item(11, 11).Value = 0.2080200501253133
$ws.Cells.Item(11, 12).Value = 0.543859649122807
$ws.Cells.Item(11, 19).Value = 0.01503759398496241
$ws.Cells.Item(12, 7).Value = 0.7048458149779736
$ws.Cells.Item(12, 10).Value = 0.2070484581497797
$ws.Cells.Item(12, 11).Value = 0.013215859030837
$ws.Cells.Item(12, 12).Value = 0.03524229074889868
$ws.Cells.Item(12, 19).Value = 0.03964757709251102
$ws.Cells.Item(13, 7).Value = 0.6511627906976745
$ws.Cells.Item(13, 10).Value = 0.3255813953488372
$ws.Cells.Item(13, 19).Value = 0.02325581395348837
$ws.Cells.Item(15, 6).Value = 0.01801801801801802
$ws.Cells.Item(15, 8).Value = 0.1396396396396396
$ws.Cells.Item(15, 9).Value = 0.09009009009009009
$ws.Cells.Item(15, 10).Value = 0.3783783783783784
$ws.Cells.Item(15, 11).Value = 0.07207207207207207
$ws.Cells.Item(15, 13).Value = 0.01351351351351351
$ws.Cells.Item(15, 15).Value = 0.05405405405405406
$ws.Cells.Item(15, 19).Value = 0.2342342342342342
$ws.Cells.Item(16, 6).Value = 0.01092896174863388
$ws.Cells.Item(16, 8).Value = 0.185792349726776
$ws.Cells.Item(16, 9).Value = 0.06557377049180328
$ws.Cells.Item(16, 10).Value = 0.3934426229508197
$ws.Cells.Item(16, 11).Value = 0.09836065573770492
$ws.Cells.Item(16, 13).Value = 0.01092896174863388
$ws.Cells.Item(16, 15).Value = 0.03278688524590164
$ws.Cells.Item(16, 19).Value = 0.2021857923497268
$ws.Cells.Item(17, 6).Value = 0.01228501228501228
$ws.Cells.Item(17, 8).Value = 0.1646191646191646
$ws.Cells.Item(17, 9).Value = 0.08353808353808354
$ws.Cells.Item(17, 10).Value = 0.3685503685503685
$ws.Cells.Item(17, 11).Value = 0.1326781326781327
$ws.Cells.Item(17, 13).Value = 0.02702702702702703
$ws.Cells.Item(17, 15).Value = 0.09582309582309582
$ws.Cells.Item(17, 19).Value = 0.1154791154791155
$ws.Cells.Item(18, 6).Value = 0.02040816326530612
$ws.Cells.Item(18, 8).Value = 0.1989795918367347
$ws.Cells.Item(18, 9).Value = 0.07653061224489796
$ws.Cells.Item(18, 10).Value = 0.413265306122449
$ws.Cells.Item(18, 11).Value = 0.07142857142857142
$ws.Cells.Item(18, 13).Value = 0.00510204081632653
$ws.Cells.Item(18, 14).Value = 0.00510204081632653
$ws.Cells.Item(18, 15).Value = 0.06122448979591837
$ws.Cells.Item(18, 19).Value = 0.1479591836734694
$ws.Cells.Item(19, 6).Value = 0.01665510062456627
$ws.Cells.Item(19, 8).Value = 0.2095766828591256
$ws.Cells.Item(19, 9).Value = 0.08119361554476058
$ws.Cells.Item(19, 10).Value = 0.3379597501734906
$ws.Cells.Item(19, 11).Value = 0.1422623178348369
$ws.Cells.Item(19, 13).Value = 0.01804302567661346
$ws.Cells.Item(19, 14).Value = 0.002775850104094379
$ws.Cells.Item(19, 15).Value = 0.06453851492019431
$ws.Cells.Item(19, 19).Value = 0.1269951422623178
